$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'55.064.53"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +7.50%  '

$ws.Range('D3').Value = "'2.465.29"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +9.30%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = "'480.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +11.43%  '

$ws.Range('D6').Value = "'139.24"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +19.31%  '

$ws.Range('D7').Value = "'0.998"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').Value = "'0.501"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +10.36%  '

$ws.Range('D9').Value = "'2.454.59"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.69%  '

$ws.Range('D10').Value = "'0.0966"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +14.49%  '

$ws.Range('D11').Value = "'5.45"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.91%  '

$ws.Range('E12').Value = '  +10.52%  '

$ws.Range('D13').Value = "'0.122"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.61%  '

$ws.Range('D14').Value = "'2.873.40"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.59%  '

$ws.Range('D15').Value = "'55.136.10"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.37%  '

$ws.Range('D16').Value = "'20.45"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +11.66%  '

$ws.Range('E17').Value = '  +19.01%  '

$ws.Range('D18').Value = "'2.447.19"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.67%  '

$ws.Range('E19').Value = '  +12.62%  '

$ws.Range('D20').Value = "'9.94"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +16.46%  '

$ws.Range('D21').Value = "'313.31"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.92%  '

$ws.Range('D22').Value = "'0.997"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.71%  '

$ws.Range('E23').Value = '  +13.28%  '

$ws.Range('D24').Value = "'57.27"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +9.95%  '

$ws.Range('E25').Value = '  +12.88%  '

$ws.Range('E26').Value = '  +0.38%  '

$ws.Range('D27').Value = "'0.162"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +16.79%  '

$ws.Range('D28').Value = "'2.548.91"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.88%  '

$ws.Range('E29').Value = '  +9.98%  '

$ws.Range('D30').Value = "'0.0₃0773"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +23.21%  '

$ws.Range('D31').Value = "'1.00"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.26%  '

$ws.Range('D32').Value = "'149.01"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.99%  '

$ws.Range('D33').Value = "'17.91"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.41%  '

$ws.Range('E34').Value = '  +13.85%  '

$ws.Range('E35').Value = '  +12.57%  '

$ws.Range('D36').Value = "'1.11"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +15.32%  '

$ws.Range('E37').Value = '  +9.92%  '

$ws.Range('E38').Value = '  +9.80%  '

$ws.Range('D39').Value = "'0.994"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.15%  '

$ws.Range('D40').Value = "'33.12"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.92%  '

$ws.Range('D41').Value = "'0.600"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.87%  '

$ws.Range('E42').Value = '  +12.33%  '

$ws.Range('E43').Value = '  +11.19%  '

$ws.Range('E44').Value = '  +14.46%  '

$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = "'10.15"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.49%  '

$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = "'256.53"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +34.61%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'4.61"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +19.32%  '

$ws.Range('D48').Value = "'0.0894"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.89%  '

$ws.Range('D49').Value = "'0.0223"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +12.25%  '

$ws.Range('D50').Value = "'1.927.05"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.53%  '

$ws.Range('D51').Value = "'16.95"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.37%  '
